$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.596.72"
$ws.Range("E2").Value = "  +4.91%  "

$ws.Range("D3").Value = "3.491.71"
$ws.Range("E3").Value = "  +3.31%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.489.93"
$ws.Range("E8").Value = "  +3.20%  "

$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.03%  "

$ws.Range("E11").Value = "  +7.18%  "

$ws.Range("E12").Value = "  +4.57%  "

$ws.Range("D13").Value = "4.096.84"
$ws.Range("E13").Value = "  +3.30%  "

$ws.Range("E14").Value = "  +0.01%  "

$ws.Range("E15").Value = "  +5.32%  "

$ws.Range("E16").Value = "  +4.33%  "

$ws.Range("D17").Value = "66.630.77"
$ws.Range("E17").Value = "  +4.87%  "

$ws.Range("D18").Value = "3.494.00"
$ws.Range("E18").Value = "  +2.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.10"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.94%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.94"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.30"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.57%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  +5.00%  "

$ws.Range("E26").Value = "  +7.53%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.69%  "

$ws.Range("E28").Value = "  +2.67%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.38%  "

$ws.Range("E31").Value = "  +6.87%  "

$ws.Range("E32").Value = "  +3.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.10%  "

$ws.Range("E35").Value = "  +0.11%  "

$ws.Range("E36").Value = "  +10.61%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.887"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.21%  "

$ws.Range("E39").Value = "  +6.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +7.71%  "

$ws.Range("E41").Value = "  +4.07%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "27.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.03%  "

$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.22%  "

$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.784.14"
$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0313"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.65%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.50"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "348.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.87%  "

$ws.Range("E50").Value = "  +7.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.89%  "
